$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 9230.0
$ws.Range("E7").Value = 480.0
$ws.Range("E8").Value = 1960.0
$ws.Range("E9").Value = 8230.0
$ws.Range("E27").Value = 730.0
$ws.Range("E29").Value = 0.0
$ws.Range("E35").Value = 0.0
$ws.Range("E36").Value = 0.0
$ws.Range("E62").Value = 155.0
$ws.Range("E84").Value = 0.0
$ws.Range("E92").Value = 0.0
$ws.Range("E104").Value = 225.0
